# "Made transition from Max->Viewer seamless."
#
# The two Scene Exporter tasks in rows 2 and 3 ("Make exporter a GUP. Build
# UI & hook data into the max files" and "Make the path from Max->Model
# Viewer seamless") are finished and are removed from the ToDo list. All
# rows below shift up by two.
#
# Two review comments are anchored on cells further down the sheet
# (originally B12 and B16); since they live below the deleted rows they
# need to move up to B10 and B14 respectively so they stay attached to the
# same logical task ("Compiled shaders" / "Error handling strategy in
# Model Compiler").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing comment text before anything moves, then remove the
# old comments (they don't automatically follow their cells when rows are
# deleted).
$shaderComment = $ws.Range("B12").Comment
$shaderCommentText = $shaderComment.Text()
$shaderComment.Delete()

$modelCompilerComment = $ws.Range("B16").Comment
$modelCompilerCommentText = $modelCompilerComment.Text()
$modelCompilerComment.Delete()

# Select and delete the two completed "Scene Exporter" rows, shifting
# everything below them up.
$ws.Rows("2:3").Select()
$ws.Rows("2:3").Delete()

# Re-attach the comments two rows higher than where they used to be, so
# they stay on the same task rows as before.
$ws.Range("B10").AddComment($shaderCommentText)
$ws.Range("B14").AddComment($modelCompilerCommentText)
